{"js": "// Update the multiplication problems (three-digit x one-digit) in the\n// practice-sheet table to the new set of values from the commit.\nconst pairs = [\n  [\"589\u00d72=\", \"291\u00d75=\"],\n  [\"615\u00d74=\", \"263\u00d76=\"],\n  [\"460\u00d79=\", \"913\u00d76=\"],\n  [\"770\u00d75=\", \"173\u00d73=\"],\n  [\"213\u00d76=\", \"417\u00d77=\"],\n  [\"668\u00d77=\", \"754\u00d75=\"],\n  [\"629\u00d76=\", \"446\u00d76=\"],\n  [\"224\u00d76=\", \"939\u00d73=\"],\n  [\"811\u00d73=\", \"562\u00d75=\"],\n  [\"160\u00d74=\", \"452\u00d77=\"],\n  [\"378\u00d78=\", \"456\u00d75=\"],\n  [\"185\u00d74=\", \"778\u00d72=\"],\n  [\"936\u00d72=\", \"301\u00d73=\"],\n  [\"665\u00d77=\", \"670\u00d72=\"],\n  [\"838\u00d72=\", \"626\u00d73=\"],\n  [\"635\u00d79=\", \"102\u00d78=\"],\n  [\"720\u00d79=\", \"664\u00d76=\"],\n  [\"751\u00d75=\", \"223\u00d73=\"],\n  [\"240\u00d78=\", \"626\u00d75=\"],\n  [\"248\u00d78=\", \"428\u00d72=\"],\n  [\"334\u00d72=\", \"954\u00d77=\"],\n  [\"970\u00d74=\", \"297\u00d72=\"],\n  [\"911\u00d72=\", \"464\u00d77=\"],\n  [\"813\u00d79=\", \"679\u00d75=\"],\n  [\"148\u00d72=\", \"913\u00d75=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems (three-digit x one-digit) in the\n# practice-sheet table to the new set of values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"589\u00d72=\", \"291\u00d75=\"),\n    @(\"615\u00d74=\", \"263\u00d76=\"),\n    @(\"460\u00d79=\", \"913\u00d76=\"),\n    @(\"770\u00d75=\", \"173\u00d73=\"),\n    @(\"213\u00d76=\", \"417\u00d77=\"),\n    @(\"668\u00d77=\", \"754\u00d75=\"),\n    @(\"629\u00d76=\", \"446\u00d76=\"),\n    @(\"224\u00d76=\", \"939\u00d73=\"),\n    @(\"811\u00d73=\", \"562\u00d75=\"),\n    @(\"160\u00d74=\", \"452\u00d77=\"),\n    @(\"378\u00d78=\", \"456\u00d75=\"),\n    @(\"185\u00d74=\", \"778\u00d72=\"),\n    @(\"936\u00d72=\", \"301\u00d73=\"),\n    @(\"665\u00d77=\", \"670\u00d72=\"),\n    @(\"838\u00d72=\", \"626\u00d73=\"),\n    @(\"635\u00d79=\", \"102\u00d78=\"),\n    @(\"720\u00d79=\", \"664\u00d76=\"),\n    @(\"751\u00d75=\", \"223\u00d73=\"),\n    @(\"240\u00d78=\", \"626\u00d75=\"),\n    @(\"248\u00d78=\", \"428\u00d72=\"),\n    @(\"334\u00d72=\", \"954\u00d77=\"),\n    @(\"970\u00d74=\", \"297\u00d72=\"),\n    @(\"911\u00d72=\", \"464\u00d77=\"),\n    @(\"813\u00d79=\", \"679\u00d75=\"),\n    @(\"148\u00d72=\", \"913\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
